$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.120.92"
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = "'1.883.77"
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = "'243.47"
$ws.Range("E5").Value = '  -2.42%  '
$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").Value = "'0.4969"
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("D8").Value = "'44.50"
$ws.Range("E8").Value = '  -2.83%  '
$ws.Range("D9").Value = "'0.2918"
$ws.Range("E9").Value = '  +1.80%  '
$ws.Range("D10").Value = "'0.06604"
$ws.Range("E10").Value = '  +0.42%  '
$ws.Range("D11").Value = "'1.879.69"
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("D12").Value = "'16.84"
$ws.Range("E12").Value = '  -1.98%  '
$ws.Range("D13").Value = "'0.07202"
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("D14").Value = "'0.6629"
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("D15").Value = "'85.60"
$ws.Range("D16").Value = "'4.836"
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").Value = "'30.102.96"
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").Value = "'0.000007860"
$ws.Range("E18").Value = '  +4.49%  '
$ws.Range("D19").Value = "'0.9997"
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").Value = "'12.77"
$ws.Range("E20").Value = '  -1.14%  '
$ws.Range("D21").Value = "'2.121.35"
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = "'0.9985"
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").Value = "'4.762"
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("D24").Value = "'5.602"
$ws.Range("E24").Value = '  +1.20%  '
$ws.Range("D25").Value = "'9.126"
$ws.Range("E25").Value = '  +0.99%  '
$ws.Range("D26").Value = "'151.50"
$ws.Range("E26").Value = '  +4.91%  '
$ws.Range("D27").Value = "'134.70"
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("D29").Value = "'1.911"
$ws.Range("E29").Value = '  -3.13%  '
$ws.Range("D30").Value = "'1.381"
$ws.Range("E30").Value = '  -0.98%  '
$ws.Range("D31").Value = "'4.159"
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").Value = "'0.08685"
$ws.Range("E32").Value = '  +0.94%  '
$ws.Range("D33").Value = "'3.938"
$ws.Range("E33").Value = '  +1.47%  '
$ws.Range("D34").Value = "'0.04995"
$ws.Range("E34").Value = '  -1.93%  '
$ws.Range("D35").Value = "'0.7086"
$ws.Range("E35").Value = '  +3.00%  '
$ws.Range("D36").Value = "'1.105"
$ws.Range("E36").Value = '  -3.31%  '
$ws.Range("D37").Value = "'2.653"
$ws.Range("E37").Value = '  -2.07%  '
$ws.Range("D38").Value = "'2.702"
$ws.Range("E38").Value = '  -2.02%  '
$ws.Range("D39").Value = "'2.192"
$ws.Range("E39").Value = '  -4.75%  '
$ws.Range("D40").Value = "'0.9347"
$ws.Range("E40").Value = '  -2.87%  '
$ws.Range("D41").Value = "'0.01648"
$ws.Range("E41").Value = '  +0.95%  '
$ws.Range("D42").Value = "'5.960"
$ws.Range("E42").Value = '  -2.37%  '
$ws.Range("D43").Value = "'0.9994"
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("E44").Value = '  -0.53%  '
$ws.Range("D45").Value = "'102.32"
$ws.Range("E45").Value = '  -1.45%  '
$ws.Range("D46").Value = "'7.486"
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("D47").Value = "'0.1255"
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").Value = "'0.05702"
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("D49").Value = "'32.39"
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("D50").Value = "'8.269"
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").Value = "'1.342"
$ws.Range("E51").Value = '  +0.28%  '
